# Generate Report for Handoff
# Update the localization status for the fb95c18c-37b4-4151-91d2-cf4ff29d38f6.md file:
# it is now "Ready for handoff" (instead of "Handed back: in sync with en-US"),
# and its latest handoff datetime is refreshed for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-09 10:43:05"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-09 10:43:10"
